# Update currency year to 2019 dollars
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Large Output Currency Unit label: "billion 2018 dollars" -> "billion 2019 dollars"
$ws.Range("A18").Value = "billion 2019 dollars"

# Medium Output Currency Unit label: "million 2018 dollars" -> "million 2019 dollars"
$ws.Range("A21").Value = "million 2019 dollars"

# Updated conversion factor (2012 dollars per 2019 dollar)
$ws.Range("A26").Value = 0.89805481563188172

# Units label for the conversion factor: "2018 dollars per 2012 dollar" -> "2019 dollars per 2012 dollar"
$ws.Range("B26").Value = "2019 dollars per 2012 dollar"

# Explanatory note referencing the conversion factor year
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2019 dollar."'
